$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.Value = "'" + $val
    $c.Style = "Normal"
}

# Price (column D) updates - forced to text to preserve exact string formatting
Set-TextValue "D2" '96.734.48'
Set-TextValue "D3" '3.615.84'
Set-TextValue "D5" '241.89'
Set-TextValue "D6" '642.29'
Set-TextValue "D10" '1.02'
Set-TextValue "D11" '3.610.18'
Set-TextValue "D12" '43.62'
Set-TextValue "D15" '4.304.52'
Set-TextValue "D16" '96.734.62'
Set-TextValue "D18" '3.616.31'
Set-TextValue "D19" '8.05'
Set-TextValue "D20" '13.12'
Set-TextValue "D21" '18.28'
Set-TextValue "D22" '0.501'
Set-TextValue "D24" '517.10'
Set-TextValue "D25" '0.0000201'
Set-TextValue "D26" '6.79'
Set-TextValue "D27" '98.46'
Set-TextValue "D28" '12.56'
Set-TextValue "D29" '3.19'
Set-TextValue "D31" '11.68'
Set-TextValue "D32" '0.999'
Set-TextValue "D35" '30.56'
Set-TextValue "D36" '0.572'
Set-TextValue "D37" '582.77'
Set-TextValue "D39" '1.48'
Set-TextValue "D49" '54.28'
Set-TextValue "D50" '8.23'

# Volume(1h) (column E) updates - plain strings (leading/trailing spaces keep them as text)
$ws.Range("E2").Value = '  +5.62%  '
$ws.Range("E3").Value = '  +8.85%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("E5").Value = '  +4.91%  '
$ws.Range("E6").Value = '  +4.61%  '
$ws.Range("E7").Value = '  +6.56%  '
$ws.Range("E8").Value = '  +5.15%  '
$ws.Range("E9").Value = '  -0.09%  '
$ws.Range("E10").Value = '  +7.61%  '
$ws.Range("E11").Value = '  +8.73%  '
$ws.Range("E12").Value = '  +3.76%  '
$ws.Range("E13").Value = '  +4.28%  '
$ws.Range("E14").Value = '  +7.79%  '
$ws.Range("E15").Value = '  +9.20%  '
$ws.Range("E16").Value = '  +5.96%  '
$ws.Range("E17").Value = '  +5.01%  '
$ws.Range("E18").Value = '  +9.08%  '
$ws.Range("E19").Value = '  -0.30%  '
$ws.Range("E20").Value = '  +21.03%  '
$ws.Range("E21").Value = '  +6.37%  '
$ws.Range("E22").Value = '  +11.93%  '
$ws.Range("E23").Value = '  +2.24%  '
$ws.Range("E24").Value = '  +5.57%  '
$ws.Range("E25").Value = '  +10.70%  '
$ws.Range("E26").Value = '  +11.35%  '
$ws.Range("E27").Value = '  +10.22%  '
$ws.Range("E28").Value = '  +6.56%  '
$ws.Range("E29").Value = '  +23.04%  '
$ws.Range("E30").Value = '  +5.91%  '
$ws.Range("E31").Value = '  +5.87%  '
$ws.Range("E32").Value = '  -0.09%  '
$ws.Range("E33").Value = '  +7.25%  '
$ws.Range("E34").Value = '  +0.50%  '
$ws.Range("E35").Value = '  +8.74%  '
$ws.Range("E36").Value = '  +9.23%  '
$ws.Range("E37").Value = '  +5.33%  '
$ws.Range("E38").Value = '  +8.45%  '
$ws.Range("E39").Value = '  +8.91%  '
$ws.Range("E40").Value = '  +3.94%  '
$ws.Range("E41").Value = '  -0.04%  '
$ws.Range("E42").Value = '  +7.12%  '
$ws.Range("E43").Value = '  +6.17%  '
$ws.Range("E44").Value = '  +5.75%  '
$ws.Range("E45").Value = '  +5.94%  '
$ws.Range("E46").Value = '  +0.30%  '
$ws.Range("E47").Value = '  +5.69%  '
$ws.Range("E48").Value = '  -3.45%  '
$ws.Range("E49").Value = '  +5.06%  '
$ws.Range("E50").Value = '  +3.92%  '
$ws.Range("E51").Value = '  +4.83%  '
